$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "EMPLOYEE_ID" column (AE),
# shifting EMPLOYEE_ID..DEPARTMENT (AE:AK) one column to the right (AF:AL).
$ws.Range("AE1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("AE1").Value = "Unnamed: 0"

# Update MANAGER_ID values (now in column AG) per the new data.
$ws.Range("AG2").Value = "O50"
$ws.Range("AG3").Value = "O50"
$ws.Range("AG4").Value = "O72"
$ws.Range("AG5").Value = "O72"
$ws.Range("AG6").Value = "O72"
$ws.Range("AG7").Value = "O72"
$ws.Range("AG8").Value = "O72"
$ws.Range("AG9").Value = "O50"

# Update PROCESS value for row 4 (now in column AK).
$ws.Range("AK4").Value = "MUTHOOT"
